$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 56 (shifts existing rows 56..149 down to 57..150)
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record
$ws.Cells.Item(56, 1).Value = 4
$ws.Cells.Item(56, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(56, 3).Value = "Los Lagos"
$ws.Cells.Item(56, 4).Value = 44477
$ws.Cells.Item(56, 5).Value = 10
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100102
$ws.Cells.Item(56, 8).Value = "Cítricos"
$ws.Cells.Item(56, 9).Value = 100102006
$ws.Cells.Item(56, 10).Value = "Pomelo"
$ws.Cells.Item(56, 11).Value = "Start Ruby"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 160
$ws.Cells.Item(56, 14).Value = 11000
$ws.Cells.Item(56, 15).Value = 12000
$ws.Cells.Item(56, 16).Value = 11500
$ws.Cells.Item(56, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(56, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 19).Value = 821
$ws.Cells.Item(56, 20).Value = 14
